# Actualizacion 10 de Mayo
# Updates the statistics sheets for the new "4APM" remedial (Rescatables) students
# and refreshes the computed totals on the 1P / 2P / Final sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Estadisticos 1P" - row for group 4APM: Blancos 13->11, Aprobados 15->17,
#    Por_Apro 53.57->60.71, Promedio 7.4->7.2
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Range("D2").Value = 11
$ws1.Range("F2").Value = 17
$ws1.Range("G2").Value = 60.71
$ws1.Range("H2").Value = 7.2

# ---------------------------------------------------------------------------
# 2) "Estadisticos 2P" - row for group 4APM: Reprobados 15->17
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Range("E2").Value = 17

# ---------------------------------------------------------------------------
# 3) "Estadisticos Final" - same update as sheet 1
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Range("D2").Value = 11
$ws3.Range("F2").Value = 17
$ws3.Range("G2").Value = 60.71
$ws3.Range("H2").Value = 7.2

# ---------------------------------------------------------------------------
# 4) "Rescatables" - insert 11 new 4APM remedial students above the existing
#    4BLCM list (which shifts from rows 2-28 down to rows 13-39).
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Rescatables")

$newRows = $ws4.Rows.Item(2).Resize(11)
$newRows.Insert()
$ws4.Range("A2:G12").ClearFormats()

$students = @(
    @(19330051920223, "ALEJO",     "VASQUEZ",   "RAUL ALEJANDRO", "FÍSICA I", "4APM", 2),
    @(19330051920225, "CESPEDES",  "CRUZ",      "ERICK MANUEL",   "FÍSICA I", "4APM", 2),
    @(19330051920229, "DIAZ",      "REYES",     "IRVING YAHIR",   "FÍSICA I", "4APM", 2),
    @(19330051920230, "ESTRADA",   "CASTAÑEDA", "BRAULIO VADIR",  "FÍSICA I", "4APM", 2),
    @(19330051920232, "GASPAR",    "TEXCAHUA",  "VIANEY",         "FÍSICA I", "4APM", 2),
    @(19330051920237, "LICEA",     "RIVERA",    "QADMIEL TAMARA", "FÍSICA I", "4APM", 2),
    @(19330051420227, "REYES",     "ANDRADE",   "ALEXANDER",      "FÍSICA I", "4APM", 2),
    @(19330051920404, "REYNOSO",   "ALCARAZ",   "RENZO JHOVANI",  "FÍSICA I", "4APM", 2),
    @(19330051920240, "RIVERA",    "GOMEZ",     "MARIA DE JESUS", "FÍSICA I", "4APM", 2),
    @(19330051920242, "SANCHEZ",   "HERNANDEZ", "ANGEL EDUARDO",  "FÍSICA I", "4APM", 2),
    @(19330051920245, "VERA",      "LLAVE",     "YESUA ISIDRO",   "FÍSICA I", "4APM", 2)
)

$r = 2
foreach ($s in $students) {
    $ws4.Cells.Item($r, 1).Value = $s[0]
    $ws4.Cells.Item($r, 2).Value = $s[1]
    $ws4.Cells.Item($r, 3).Value = $s[2]
    $ws4.Cells.Item($r, 4).Value = $s[3]
    $ws4.Cells.Item($r, 5).Value = $s[4]
    $ws4.Cells.Item($r, 6).Value = $s[5]
    $ws4.Cells.Item($r, 7).Value = $s[6]
    $r = $r + 1
}
